$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells are formatted as Text so the literal string values
# (e.g. "277.73", "6.40%") are preserved exactly as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.40%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.31%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.811"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.14%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06248"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.50%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.909"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.32%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.273"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.37%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8806"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.44%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9420"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.87%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1450"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.13%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05254"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.28%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07281"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.88%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03160"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.00%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09054"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.02%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001550"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.69%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006273"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.86%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005824"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.62%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.22%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "6.51%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3094"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1293"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.34%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.849"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.14%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04327"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.79%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.10%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004263"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.03%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001691"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.11%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04033"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.06%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006384"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "54.40%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.54%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.45%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01221"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.16%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005086"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.48%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.02%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.375"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "838.15%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.02%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"

